$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update B5 / C5 with new, more specific text for BMP180 (adds temperature sensing info)
$ws.Range("B5").Value = "Digitaler Luftdrucksensor und Temperatursensor"
$ws.Range("C5").Value = "Misst Luftdruck`nAuch das Messen von Temperatur ist möglich"

# Move the active selection on the bottom-right pane to C6 (cosmetic, matches author's cursor position)
$ws.Range("C6").Select()
